# Update the "想去人数" (number of people wanting to go) figures (column F)
# for rows 3, 4, 6 and 7 on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2037
    $ws.Range("F4").Value = 247
    $ws.Range("F6").Value = 6365
    $ws.Range("F7").Value = 242
}
